$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 3.479888236813032
$ws.Range("E2").Value = 0.8548164319983614
$ws.Range("H2").Value = 3.341511580334618
$ws.Range("I2").Value = 0.7582507812126756
